# Journal.xlsx update - "Mise a jour du journal de bord"
# Fill in the next two journal entries (Github desktop learning + setup)
# in the first empty slots of the two side-by-side tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Right-hand table: row 12 (I12:N12) was the first empty row ---
# Match the date number format already used by the table (copy from the row above)
$ws.Range("J11").Copy()
$ws.Range("J12").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("I12").Value = "ICT-431"
$ws.Range("J12").Value = 43558
$ws.Range("K12").Value = "70m"
$ws.Range("L12").Value = "Théorie"
$ws.Range("M12").Value = "Apprendre comment utiliser Github desktop"

# --- Left-hand table: row 18 (B18:G18) was the first empty row ---
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B18").Value = "ICT-431"
$ws.Range("C18").Value = 43558
$ws.Range("D18").Value = "20m"
$ws.Range("E18").Value = "Documentation"
$ws.Range("F18").Value = "Installer Github et passer la documentation sur github"

# Row 18 grows a bit taller to match the new wrapped comment text
$ws.Rows.Item(18).RowHeight = 39

# Widen column N slightly so the new comment text fits better
$ws.Columns.Item(14).ColumnWidth = 29.7

# Move the active selection to where the user ended up after the edit
$ws.Range("E23").Select()
